# lower-cased username and password on login
# (two new job-log rows were appended to the JOBS sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 26 ----
$ws.Range("A26").Value = 71278
$ws.Range("B26").Value = "06DF3136AA3600"
$ws.Range("C26").Value = "3600J03031"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "460"
$ws.Range("E26").Value = "?"
$ws.Range("F26").Value = "GOOD"
$ws.Range("G26").Value = "NA"
$ws.Range("H26").Value = "NO"
$ws.Range("I26").NumberFormat = "@"
$ws.Range("I26").Value = ""
$ws.Range("J26").Value = "ravi"
$ws.Range("K26").NumberFormat = "@"
$ws.Range("K26").Value = "8/31/2022"
$ws.Range("L26").Value = $false
$ws.Range("M26").Value = "N/A"
$ws.Range("N26").Value = "N/A"
$ws.Range("O26").Value = "NO"

# ---- Row 27 ----
$ws.Range("A27").Value = 71279
$ws.Range("B27").Value = "JAS0N"
$ws.Range("C27").Value = "DRFDWR"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "460"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2"
$ws.Range("F27").Value = "GOOD"
$ws.Range("G27").Value = "N/A23345"
$ws.Range("H27").Value = "NO"
$ws.Range("I27").Value = "Jgyjgygujguj"
$ws.Range("J27").Value = "west"
$ws.Range("K27").NumberFormat = "@"
$ws.Range("K27").Value = "8/31/2022"
$ws.Range("L27").Value = $false
$ws.Range("M27").Value = "N/A"
$ws.Range("N27").Value = "N/A"
$ws.Range("O27").Value = "?"

# Keep the "number stored as text" warning suppressed over the full data
# range now that it has grown to include the two new rows.
$errors = $ws.Range("A1:O27").Errors
$errors.Item(9).Ignore = $true
